# Add 5 new suffix rows (more than 1000 dataset) to the suffix_set sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weight/suffix pairs appended right after the existing data (rows 100-103),
# following the exact same layout/formatting already used by the sheet's
# default column styles (column A = "weight", column B = "suffix").
$newRows = @(
    @(3, "টারই"),
    @(3, "টুকু"),
    @(3, "েছিলেন"),
    @(3, "েক"),
    @(3, "েই")
)

$startRow = $ws.UsedRange.Rows.Count + 1
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

$endRow = $startRow + $newRows.Count - 1

# Match the selection left behind in the sheet view: active cell on the
# second newly-added row, with the remaining new rows selected alongside it.
$ws.Range("A105:A108").Select() | Out-Null

# Best effort: keep the workbook window sized the same way the authored
# workbook ended up (not all runtimes persist this window chrome setting).
$win = $wb.Windows.Item(1)
$win.Width = 11520
$win.Height = 9072
